# Updated cryptos list - apply new Price and Volume(1h) values to rows 2-51
function Set-CellText($Worksheet, $Row, $Col, $Text) {
    $c = $Worksheet.Cells.Item($Row, $Col)
    # Prefix with an apostrophe so Excel treats numeric-looking strings
    # (e.g. "0.999", "56.50") as text rather than auto-converting them
    # to numbers; ClearFormats() then strips the resulting quote-prefix
    # text format so no stray style/number-format is left behind.
    $c.Value = "'" + $Text
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 2 4 '64.009.83'
Set-CellText $ws 2 5 '  -5.39%  '
Set-CellText $ws 3 4 '3.283.49'
Set-CellText $ws 3 5 '  -6.72%  '
Set-CellText $ws 4 4 '0.999'
Set-CellText $ws 4 5 '  -0.17%  '
Set-CellText $ws 5 4 '522.34'
Set-CellText $ws 5 5 '  -6.12%  '
Set-CellText $ws 6 4 '174.12'
Set-CellText $ws 6 5 '  -14.27%  '
Set-CellText $ws 7 4 '0.602'
Set-CellText $ws 7 5 '  -1.15%  '
Set-CellText $ws 8 4 '3.274.72'
Set-CellText $ws 8 5 '  -6.70%  '
Set-CellText $ws 9 5 '  -0.08%  '
Set-CellText $ws 10 4 '0.605'
Set-CellText $ws 10 5 '  -8.00%  '
Set-CellText $ws 11 4 '56.50'
Set-CellText $ws 11 5 '  -12.18%  '
Set-CellText $ws 12 4 '0.133'
Set-CellText $ws 12 5 '  -8.12%  '
Set-CellText $ws 13 4 '0.0000257'
Set-CellText $ws 13 5 '  -5.73%  '
Set-CellText $ws 14 4 '9.02'
Set-CellText $ws 14 5 '  -9.04%  '
Set-CellText $ws 15 4 '3.798.95'
Set-CellText $ws 15 5 '  -6.82%  '
Set-CellText $ws 16 4 '3.271.89'
Set-CellText $ws 16 5 '  -6.94%  '
Set-CellText $ws 17 4 '0.116'
Set-CellText $ws 17 5 '  -6.40%  '
Set-CellText $ws 18 4 '63.815.61'
Set-CellText $ws 18 5 '  -5.29%  '
Set-CellText $ws 19 4 '17.31'
Set-CellText $ws 19 5 '  -6.72%  '
Set-CellText $ws 20 4 '11.02'
Set-CellText $ws 20 5 '  -6.99%  '
Set-CellText $ws 21 4 '0.950'
Set-CellText $ws 21 5 '  -7.78%  '
Set-CellText $ws 22 4 '371.81'
Set-CellText $ws 22 5 '  -5.47%  '
Set-CellText $ws 23 4 '3.74'
Set-CellText $ws 23 5 '  -6.88%  '
Set-CellText $ws 24 4 '80.04'
Set-CellText $ws 24 5 '  -4.00%  '
Set-CellText $ws 25 4 '10.97'
Set-CellText $ws 25 5 '  -9.99%  '
Set-CellText $ws 26 4 '3.86'
Set-CellText $ws 26 5 '  -1.60%  '
Set-CellText $ws 27 5 '  -1.95%  '
Set-CellText $ws 28 4 '2.65'
Set-CellText $ws 28 5 '  -6.67%  '
Set-CellText $ws 29 4 '11.31'
Set-CellText $ws 29 5 '  -7.67%  '
Set-CellText $ws 30 4 '8.27'
Set-CellText $ws 30 5 '  -6.92%  '
Set-CellText $ws 31 4 '28.63'
Set-CellText $ws 31 5 '  -8.03%  '
Set-CellText $ws 32 4 '639.13'
Set-CellText $ws 32 5 '  -10.65%  '
Set-CellText $ws 33 4 '6.62'
Set-CellText $ws 33 5 '  -6.75%  '
Set-CellText $ws 34 4 '11.19'
Set-CellText $ws 34 5 '  -5.10%  '
Set-CellText $ws 35 4 '0.105'
Set-CellText $ws 35 5 '  -6.38%  '
Set-CellText $ws 36 4 '58.69'
Set-CellText $ws 36 5 '  -8.36%  '
Set-CellText $ws 38 4 '36.39'
Set-CellText $ws 38 5 '  -5.99%  '
Set-CellText $ws 39 4 '0.383'
Set-CellText $ws 39 5 '  -3.99%  '
Set-CellText $ws 40 4 '0.997'
Set-CellText $ws 40 5 '  -0.05%  '
Set-CellText $ws 41 4 '0.0₃0695'
Set-CellText $ws 41 5 '  +0.78%  '
Set-CellText $ws 42 5 '  -6.44%  '
Set-CellText $ws 43 4 '2.903.46'
Set-CellText $ws 43 5 '  -5.13%  '
Set-CellText $ws 44 4 '2.44'
Set-CellText $ws 44 5 '  -6.34%  '
Set-CellText $ws 45 4 '2.66'
Set-CellText $ws 45 5 '  -11.65%  '
Set-CellText $ws 46 5 '  -4.97%  '
Set-CellText $ws 47 4 '0.0394'
Set-CellText $ws 47 5 '  -3.55%  '
Set-CellText $ws 48 5 '  +5.84%  '
Set-CellText $ws 49 4 '0.125'
Set-CellText $ws 49 5 '  -1.87%  '
Set-CellText $ws 50 4 '2.74'
Set-CellText $ws 50 5 '  +4.13%  '
Set-CellText $ws 51 4 '134.32'
Set-CellText $ws 51 5 '  -3.07%  '
